$wb = $excel.ActiveWorkbook

# New material row (index 16) to be appended to every data sheet:
# "Finemet FT-3M .018mm"
$matName = "Finemet FT-3M .018mm"

# Sheet "Freq": F1/F2 values for the new material
$ws = $wb.Worksheets.Item("Freq")
$ws.Range("A3").Value = 16
$ws.Range("B3").Value = $matName
$ws.Range("C3").Value = 20000
$ws.Range("D3").Value = 20000
$ws.Range("E3").Value = 10000
$ws.Range("F3").Value = 10000
$ws.Activate() | Out-Null
$ws.Range("B11").Select() | Out-Null

# Sheet "Bfield": B field values for the new material
$ws = $wb.Worksheets.Item("Bfield")
$ws.Range("A3").Value = 16
$ws.Range("B3").Value = $matName
$ws.Range("C3").Value = 0.3
$ws.Range("D3").Value = 0.6
$ws.Range("E3").Value = 0.3
$ws.Range("F3").Value = 0.6
$ws.Activate() | Out-Null
$ws.Range("A3").Select() | Out-Null

# Sheet "BSAT": saturation flux density for the new material
$ws = $wb.Worksheets.Item("BSAT")
$ws.Range("A3").Value = 16
$ws.Range("B3").Value = $matName
$ws.Range("C3").Value = 1.23
$ws.Activate() | Out-Null
$ws.Range("C7").Select() | Out-Null

# Sheet "Ploss": core loss, computed with formulas for the new material
$ws = $wb.Worksheets.Item("Ploss")
$ws.Range("A3").Value = 16
$ws.Range("B3").Value = $matName
$ws.Range("C3").Formula = "=7.3*4.482"
$ws.Range("D3").Formula = "=7.3*18.397"
$ws.Range("E3").Formula = "=7.3*3.762"
$ws.Range("F3").Formula = "=7.3*16.398"
$ws.Activate() | Out-Null
$ws.Range("D22").Select() | Out-Null

# Sheet "MU": permeability for the new material
$ws = $wb.Worksheets.Item("MU")
$ws.Range("A3").Value = 16
$ws.Range("B3").Value = $matName
$ws.Range("C3").Value = 70000
$ws.Activate() | Out-Null
$ws.Range("C9").Select() | Out-Null

# Sheet "Density": density for the new material
$ws = $wb.Worksheets.Item("Density")
$ws.Range("A3").Value = 16
$ws.Range("B3").Value = $matName
$ws.Range("C3").Value = 7.7
$ws.Activate() | Out-Null
$ws.Range("D9").Select() | Out-Null
